$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must remain stored as literal text
# (matches the source data which stores these as plain text/inline strings).
# Using NumberFormat "@" forces text entry, then resetting the Style back to
# "Normal" avoids leaving a stray per-cell format applied.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7405"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3152"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07196"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08381"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7506"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.428"
$ws.Range("D13").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.075"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "246.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007847"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9984"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.009"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1555"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.262"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.037"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.497"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.601"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.538"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.279"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05316"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.238"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7541"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9998"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.691"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01962"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.756"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4501"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.058"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8565"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.630"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.858"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.470"
$ws.Range("D50").Style = "Normal"

# Remaining cells (percentages, names, URLs, and D-values that already
# contain non-numeric punctuation) can be set directly as text.
$ws.Range("D2").Value = "29.908.95"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.875.46"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -4.41%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  -3.72%  "
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").Value = "1.889.27"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "29.908.19"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "2.126.29"
$ws.Range("E22").Value = "  -4.12%  "
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -4.99%  "
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +4.15%  "
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  +4.19%  "
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").Value = "1.108.77"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  -3.22%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.024.30"
$ws.Range("E51").Value = "  -4.20%  "
